$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "MOHD ARIFFIN BIN AHMAD"
$ws.Range("C2").Value = "+6013- 8693559"
$ws.Range("D2").Value = "ariffinahmad1994@gmail.com"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "['Malaysia', 'Sabah', 'Tawau']"
$ws.Range("H2").Value = "[{'field_of_study': 'Geology', 'level': ""Bachelor's Degree"", 'cgpa': 'N/A', 'university': 'Universiti Malaysia Sabah', 'start_date': '2013-09-01', 'year_of_graduation': '2016'}]"
$ws.Range("I2").Value = "['PTW Level 1 (INSTEP)', 'PTW Level 2 (INSTEP)', 'OGSP (NIOSH)', 'First Aid Training (Red Crescent, Miri)', 'Basic Fire Fighting Training (Bomba, Lintas)', 'CIDB']"
$ws.Range("J2").Value = "['N/A']"
$ws.Range("K2").Value = "['AutoCAD 20 21', 'CorelDRAW X7', 'ArcGIS (Basic)', 'Surfer', 'GCDkit', 'Ilwis', 'Microsoft Office']"
$ws.Range("L2").Value = "[]"
$ws.Range("M2").Value = "[{'job_title': 'Geologist', 'job_company': 'SATOK BRIDGE AREA CONSULTING ENGINEERS SDN. BHD.', 'Industries': 'Engineering', 'start_date': '2020-11-01', 'end_date': '2024-05-15'}, {'job_title': 'Geologist', 'job_company': 'SATOK BRIDGE AREA CONSULTING ENGINEERS SDN. BHD.', 'Industries': 'Engineering', 'start_date': '2020-02-01', 'end_date': '2020-10-01'}, {'job_title': 'Project Coordinator', 'job_company': 'MAJUMEC BINA SDN. BHD.', 'Industries': 'Construction', 'start_date': '2018-05-01', 'end_date': '2020-01-01'}, {'job_title': 'Engineering Geologist', 'job_company': 'MAJUMEC BINA SDN. BHD.', 'Industries': 'Construction', 'start_date': '2017-08-01', 'end_date': '2018-04-01'}, {'job_title': 'Junior Geologist', 'job_company': 'GEOSPEC SDN. BHD.', 'Industries': 'Geology', 'start_date': '2016-07-01', 'end_date': '2016-09-01'}]"

# Row 3
$ws.Range("B3").Value = "DARIN NUR ADLEEN BINTI ROSALIZI"
$ws.Range("C3").Value = "+6019 353 -4359"
$ws.Range("D3").Value = "darinrosalizi@gmail.com"
$ws.Range("G3").Value = "['Malaysia', 'Selangor', 'Sungai Buloh']"
$ws.Range("H3").Value = "[{'field_of_study': 'Intelligent Systems Engineering', 'level': 'Bachelor of Information Systems (Honors)', 'cgpa': '3.67', 'university': 'UiTM Shah Alam', 'start_date': '2020-03', 'year_of_graduation': '2022'}, {'field_of_study': 'Computer Science', 'level': 'Diploma', 'cgpa': '3.67', 'university': 'UiTM Segamat', 'start_date': '2017-05', 'year_of_graduation': '2020'}]"
$ws.Range("J3").Value = "['Requirement Analysis', 'Technical Documentation', 'Agile Methodologies', 'Project Management', 'Process Flow Diagrams']"
$ws.Range("L3").Value = "['English', 'Malay']"
$ws.Range("M3").Value = "[{'job_title': 'System Analyst', 'job_company': 'Artanis Cloud Sdn Bhd', 'Industries': 'Information Technology', 'start_date': '2022-03', 'end_date': '2022-08', 'job_location': 'Shah Alam, Selangor'}, {'job_title': 'System Analyst', 'job_company': 'Artanis Cloud Sdn Bhd', 'Industries': 'Information Technology', 'start_date': '2022-09', 'end_date': '2023-01', 'job_location': 'Shah Alam, Selangor'}, {'job_title': 'Sales Associate', 'job_company': 'TS Concept Store', 'Industries': 'Retail', 'start_date': '2017-02', 'end_date': '2017-04', 'job_location': 'Sungai Buloh, Selangor'}, {'job_title': 'System Developer', 'job_company': 'SIRIM Berhad', 'Industries': 'Research and Development', 'start_date': '2019-07', 'end_date': '2019-09', 'job_location': 'Shah Alam, Selangor'}]"
